$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 124.75
$ws.Range("I2").Value = 109.6
$ws.Range("K2").Value = 109.6
$ws.Range("M2").Value = 3.400000000000006
$ws.Range("H5").Value = 135.33333
$ws.Range("I5").Value = 197
$ws.Range("J5").Value = 12
$ws.Range("K5").Value = 197
$ws.Range("L5").Value = 12
$ws.Range("M5").Value = -82
$ws.Range("N5").Value = -242
$ws.Range("H94").Value = 1133.1666
$ws.Range("I94").Value = 1133.1666
$ws.Range("K94").Value = 1133.1666
$ws.Range("M94").Value = -682.1666
$ws.Range("H132").Value = 7367.2573
$ws.Range("I132").Value = 8196.612999999999
$ws.Range("J132").Value = 939.75
$ws.Range("K132").Value = 24589.839
$ws.Range("L132").Value = 2819.25
$ws.Range("M132").Value = -22059.839
$ws.Range("N132").Value = -7879.25
$ws.Range("H137").Value = 4552624.5
$ws.Range("I137").Value = 7145071.5
$ws.Range("J137").Value = 15843
$ws.Range("K137").Value = 21435214.5
$ws.Range("L137").Value = 47529
$ws.Range("M137").Value = -21432664.5
$ws.Range("N137").Value = -52629
$ws.Range("H138").Value = 6017.787
$ws.Range("I138").Value = 6368.615
$ws.Range("J138").Value = 5922.771
$ws.Range("K138").Value = 19105.845
$ws.Range("L138").Value = 17768.313
$ws.Range("M138").Value = -13965.845
$ws.Range("N138").Value = -28048.313
$ws.Range("H141").Value = 4981.864
$ws.Range("I141").Value = 2107.0667
$ws.Range("K141").Value = 6321.2001
$ws.Range("M141").Value = -1141.2001

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 772081.5600000001
$ws.Range("I32").Value = 894117.75
$ws.Range("J32").Value = 12745.444
$ws.Range("K32").Value = 894117.75
$ws.Range("L32").Value = 12745.444
$ws.Range("M32").Value = -893830.75
$ws.Range("N32").Value = -13319.444
$ws.Range("H74").Value = 2185181.8
$ws.Range("I74").Value = 3111474.2
$ws.Range("K74").Value = 3111474.2
$ws.Range("M74").Value = -3110600.2
$ws.Range("H77").Value = 2185181.8
$ws.Range("I77").Value = 3111474.2
$ws.Range("K77").Value = 15557371
$ws.Range("M77").Value = -15553003
$ws.Range("H110").Value = 1599.75
$ws.Range("I110").Value = 1466.3334
$ws.Range("K110").Value = 1466.3334
$ws.Range("M110").Value = 578.6666
$ws.Range("H133").Value = 90000
$ws.Range("I133").Value = 90000
$ws.Range("K133").Value = 90000
$ws.Range("M133").Value = -87470

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H99").Value = 6088.56
$ws.Range("J99").Value = 4633.3335
$ws.Range("L99").Value = 4633.3335
$ws.Range("N99").Value = -7629.3335
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H105").Value = 5070.44
$ws.Range("I105").Value = 5244.2085
$ws.Range("K105").Value = 5244.2085
$ws.Range("M105").Value = -3497.2085
$ws.Range("H122").Value = 24688
$ws.Range("I122").Value = 904.75
$ws.Range("J122").Value = 38278.43
$ws.Range("K122").Value = 2714.25
$ws.Range("L122").Value = 114835.29
$ws.Range("M122").Value = -264.25
$ws.Range("N122").Value = -119735.29

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 13528642
$ws.Range("I4").Value = 18615736
$ws.Range("K4").Value = 55847208
$ws.Range("M4").Value = -55847096
$ws.Range("H68").Value = 4849.327
$ws.Range("J68").Value = 4895.5283
$ws.Range("L68").Value = 14686.5849
$ws.Range("N68").Value = -16308.5849
$ws.Range("H71").Value = 4849.327
$ws.Range("J71").Value = 4895.5283
$ws.Range("L71").Value = 44059.7547
$ws.Range("N71").Value = -52171.7547
$ws.Range("H107").Value = 5149.9
$ws.Range("I107").Value = 3516.5
$ws.Range("K107").Value = 10549.5
$ws.Range("M107").Value = -8629.5
$ws.Range("H108").Value = 4000
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 4000
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 12000
$ws.Range("N108").Value = -17760
$ws.Range("M108").ClearContents()
$ws.Range("H109").Value = 5456.5625
$ws.Range("I109").Value = 634.3333
$ws.Range("J109").Value = 8349.9
$ws.Range("K109").Value = 1902.9999
$ws.Range("L109").Value = 25049.7
$ws.Range("M109").Value = -862.9999
$ws.Range("N109").Value = -27129.7
$ws.Range("H122").Value = 769110.7
$ws.Range("I122").Value = 1344725.5
$ws.Range("K122").Value = 12102529.5
$ws.Range("M122").Value = -12100079.5
$ws.Range("H138").Value = 17251.1
$ws.Range("I138").Value = 20252.588
$ws.Range("K138").Value = 60757.764
$ws.Range("M138").Value = -55617.764

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 7002
$ws.Range("J10").Value = 7002
$ws.Range("L10").Value = 7002
$ws.Range("N10").Value = -7340
$ws.Range("H107").Value = 1038.8
$ws.Range("I107").Value = 1046.8182
$ws.Range("K107").Value = 1046.8182
$ws.Range("M107").Value = 873.1818000000001
$ws.Range("H126").Value = 2881.4285
$ws.Range("I126").Value = 2355.25
$ws.Range("J126").Value = 3583
$ws.Range("K126").Value = 7065.75
$ws.Range("L126").Value = 10749
$ws.Range("M126").Value = -4595.75
$ws.Range("N126").Value = -15689

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2979943.8
$ws.Range("I132").Value = 4389333.5
$ws.Range("J132").Value = 4566.1113
$ws.Range("K132").Value = 13168000.5
$ws.Range("L132").Value = 13698.3339
$ws.Range("M132").Value = -13165470.5
$ws.Range("N132").Value = -18758.3339
$ws.Range("H136").Value = 15627253
$ws.Range("I136").Value = 10418921
$ws.Range("J136").Value = 31252250
$ws.Range("K136").Value = 31256763
$ws.Range("L136").Value = 93756750
$ws.Range("M136").Value = -31254213
$ws.Range("N136").Value = -93761850

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H126").Value = 3389.3076
$ws.Range("I126").Value = 3673.889
$ws.Range("J126").Value = 2749
$ws.Range("K126").Value = 11021.667
$ws.Range("L126").Value = 8247
$ws.Range("M126").Value = -8551.667000000001
$ws.Range("N126").Value = -13187
$ws.Range("H132").Value = 6669404.5
$ws.Range("I132").Value = 9261820
$ws.Range("K132").Value = 27785460
$ws.Range("M132").Value = -27782930

Write-Host "All updates applied."